$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = "18_hazards_to_humans_and_domestic_animals"
$ws.Range("F5").Value = "application instructions || env warning - species || pollinator"
$ws.Range("F6").Value = "32_physical_and_chemical_hazards"
$ws.Range("F8").Value = "135_product_information"
$ws.Range("F9").Value = "application instructions"
$ws.Range("F10").Value = "mixing"
$ws.Range("F11").Value = "application instructions"
$ws.Range("F15").Value = "application instructions"
$ws.Range("F16").Value = "application instructions"
$ws.Range("F17").Value = "application instructions"
$ws.Range("F19").Value = "application instructions"
$ws.Range("F20").Value = "application instructions"
$ws.Range("F45").Value = "application instructions"
$ws.Range("F46").Value = "application instructions"
$ws.Range("F47").Value = "application instructions"
$ws.Range("F50").Value = "irrigation"
